$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# "Title ID 1" (G1) / "Title ID 2" (M1) become rich-text "Film Code N
# (internal reference, same as other excel import)" headers, bold + white,
# with "internal" additionally underlined — using film code instead of
# titleId.
# ---------------------------------------------------------------------------

function Set-FilmCodeHeader($cellRef, $n) {
    $nl = [char]10
    $text = "Film Code " + $n + " " + $nl + "(internal reference, same as other excel import)"
    $cell = $ws.Range($cellRef)
    $cell.Value = $text

    # "Film Code N <newline>(" -> bold, white, Arial 10, no underline
    $run1 = $cell.Characters(1, 14)
    $run1.Font.Name = "Arial"
    $run1.Font.Size = 10
    $run1.Font.Bold = $true
    $run1.Font.Underline = $false
    $run1.Font.Color = 16777215

    # "internal" -> bold, white, Arial 10, underlined
    $run2 = $cell.Characters(15, 8)
    $run2.Font.Name = "Arial"
    $run2.Font.Size = 10
    $run2.Font.Bold = $true
    $run2.Font.Underline = $true
    $run2.Font.Color = 16777215

    # " reference, same as other excel import)" -> bold, white, Arial 10, no underline
    $run3 = $cell.Characters(23, 39)
    $run3.Font.Name = "Arial"
    $run3.Font.Size = 10
    $run3.Font.Bold = $true
    $run3.Font.Underline = $false
    $run3.Font.Color = 16777215

    # Whole-cell font stays bold / white / Arial 10, no underline (matches the
    # surrounding header row styling, fill colour is untouched).
    $cell.Font.Name = "Arial"
    $cell.Font.Size = 10
    $cell.Font.Bold = $true
    $cell.Font.Underline = $false
    $cell.Font.Color = 16777215
}

Set-FilmCodeHeader "G1" 1
Set-FilmCodeHeader "M1" 2

# ---------------------------------------------------------------------------
# Header row 1 grows taller to fit the two-line film-code headers.
# ---------------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 35.25

# ---------------------------------------------------------------------------
# Selection / scroll position: now parked on G1 (was B31).
# ---------------------------------------------------------------------------
$ws.Range("G1").Select()
